$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.903.18"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.702.32"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.95%  "

$ws.Range("D5").Value = "315.70"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D7").Value = "0.4029"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").Value = "1.004"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").Value = "53.63"
$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("D11").Value = "1.468"
$ws.Range("E11").Value = "  -3.63%  "

$ws.Range("D12").Value = "0.08808"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("D13").Value = "25.95"
$ws.Range("E13").Value = "  +5.23%  "

$ws.Range("D14").Value = "7.484"
$ws.Range("E14").Value = "  -2.43%  "

$ws.Range("D15").Value = "8.033"
$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").Value = "0.00001348"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("D17").Value = "1.778.84"
$ws.Range("E17").Value = "  +3.79%  "

$ws.Range("D18").Value = "96.43"
$ws.Range("E18").Value = "  -3.96%  "

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "21.11"
$ws.Range("E20").Value = "  +5.25%  "

$ws.Range("D21").Value = "7.248"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("D23").Value = "14.51"
$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").Value = "24.895.70"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  -1.18%  "

$ws.Range("D26").Value = "2.897"
$ws.Range("E26").Value = "  -5.11%  "

$ws.Range("D27").Value = "6.572"
$ws.Range("E27").Value = "  +26.25%  "

$ws.Range("D28").Value = "23.03"
$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("D29").Value = "165.47"
$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").Value = "145.08"
$ws.Range("E30").Value = "  +4.25%  "

$ws.Range("D31").Value = "8.242"
$ws.Range("E31").Value = "  -5.11%  "

$ws.Range("D32").Value = "1.931.43"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").Value = "2.259"
$ws.Range("E33").Value = "  +13.29%  "

$ws.Range("D34").Value = "0.08811"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").Value = "0.03212"
$ws.Range("E35").Value = "  +8.75%  "

$ws.Range("D36").Value = "7.311"
$ws.Range("E36").Value = "  -4.46%  "

$ws.Range("D37").Value = "1.021"
$ws.Range("E37").Value = "  -3.57%  "

$ws.Range("E38").Value = "  +3.16%  "

$ws.Range("D39").Value = "0.8429"
$ws.Range("E39").Value = "  +5.79%  "

$ws.Range("D40").Value = "10.84"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").Value = "0.09364"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("D42").Value = "14.08"
$ws.Range("E42").Value = "  -3.03%  "

$ws.Range("D43").Value = "1.472"
$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "17.82"
$ws.Range("E44").Value = "  +7.46%  "

$ws.Range("D45").Value = "2.718"
$ws.Range("E45").Value = "  +3.54%  "

$ws.Range("D46").Value = "0.7444"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  -0.65%  "

$ws.Range("D48").Value = "1.392"
$ws.Range("E48").Value = "  +4.40%  "

$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("D50").Value = "142.31"
$ws.Range("E50").Value = "  +1.62%  "

$ws.Range("D51").Value = "0.08337"
$ws.Range("E51").Value = "  +3.42%  "
